$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8617596626281738
$ws.Range("B1").Value = 1.208746194839478
$ws.Range("C1").Value = 1.707288146018982
$ws.Range("D1").Value = 5.116880416870117
$ws.Range("E1").Value = 2.140952587127686
